# Update the "Corr/total marks" figures on the concise marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total correct answers (row "Marking", column "Right")
$ws.Range("B11").Value = 5

# Total score (row "Total", column "Right")
$ws.Range("B12").Value = 105

# Correct/total marks summary text
$ws.Range("E12").Value = "105/140"
